# Update the Fitness (column C) values in Sheet1 to reflect the new run data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 2-39 (Generation 0-37): Fitness 7622 -> 7312
$ws.Range("C2:C39").Value = 7312

# Rows 40-118 (Generation 38-116): Fitness 7622 -> 7310
$ws.Range("C40:C118").Value = 7310

# Rows 119-252 (Generation 117-250): Fitness 7622 -> 7293
$ws.Range("C119:C252").Value = 7293
